# chore: update Sheets via scheduled runner
# Refresh market-price derived columns (currentAveragePrice / currentAveragePriceNQ /
# currentAveragePriceHQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
# for the affected leve rows across each crafting-class sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 12988175
$ws.Range("I112").Value = 333.33334
$ws.Range("J112").Value = 13938505
$ws.Range("K112").Value = 1000.00002
$ws.Range("L112").Value = 41815515
$ws.Range("M112").Value = 107.9999799999999
$ws.Range("N112").Value = -41817731

$ws.Range("H138").Value = 1855.54
$ws.Range("I138").Value = 840.52
$ws.Range("J138").Value = 2870.56
$ws.Range("K138").Value = 2521.56
$ws.Range("L138").Value = 8611.68
$ws.Range("M138").Value = 2618.44
$ws.Range("N138").Value = -18891.68

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H51").Value = 29166.666
$ws.Range("J51").Value = 29166.666
$ws.Range("L51").Value = 29166.666
$ws.Range("N51").Value = -30678.666

$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws.Range("H61").Value = 222144.38
$ws.Range("I61").Value = 5281.2583
$ws.Range("K61").Value = 5281.2583
$ws.Range("M61").Value = -5069.2583

$ws.Range("H132").Value = 1540927.8
$ws.Range("I132").Value = 1579.6305
$ws.Range("J132").Value = 5267770.5
$ws.Range("K132").Value = 4738.8915
$ws.Range("L132").Value = 15803311.5
$ws.Range("M132").Value = -2208.8915
$ws.Range("N132").Value = -15808371.5

$ws.Range("H136").Value = 222144.38
$ws.Range("I136").Value = 5281.2583
$ws.Range("K136").Value = 15843.7749
$ws.Range("M136").Value = -13293.7749

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2040.4286
$ws.Range("I105").Value = 1856.6
$ws.Range("K105").Value = 1856.6
$ws.Range("M105").Value = -109.5999999999999

$ws.Range("H132").Value = 50680
$ws.Range("J132").Value = 50680
$ws.Range("L132").Value = 50680
$ws.Range("N132").Value = -60800

$ws.Range("H134").Value = 18388.879
$ws.Range("I134").Value = 3445.7017
$ws.Range("J134").Value = 113029
$ws.Range("K134").Value = 10337.1051
$ws.Range("L134").Value = 339087
$ws.Range("M134").Value = -7802.105100000001
$ws.Range("N134").Value = -344157

$ws.Range("H140").Value = 52910
$ws.Range("J140").Value = 52910
$ws.Range("L140").Value = 52910
$ws.Range("N140").Value = -63270

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 63401.2
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H16").Value = 2052.5833
$ws.Range("I16").Value = 1582
$ws.Range("J16").Value = 2523.1667
$ws.Range("K16").Value = 1582
$ws.Range("L16").Value = 2523.1667
$ws.Range("M16").Value = -1295
$ws.Range("N16").Value = -3097.1667

$ws.Range("H31").Value = 295429.4
$ws.Range("I31").Value = 1663.6786
$ws.Range("K31").Value = 1663.6786
$ws.Range("M31").Value = -1368.6786

$ws.Range("H34").Value = 295429.4
$ws.Range("I34").Value = 1663.6786
$ws.Range("K34").Value = 1663.6786
$ws.Range("M34").Value = -1461.6786

$ws.Range("H113").Value = 2052.5833
$ws.Range("I113").Value = 1582
$ws.Range("J113").Value = 2523.1667
$ws.Range("K113").Value = 1582
$ws.Range("L113").Value = 2523.1667
$ws.Range("M113").Value = 588
$ws.Range("N113").Value = -6863.1667

$ws.Range("H122").Value = 1960.4359
$ws.Range("I122").Value = 1565.875
$ws.Range("J122").Value = 2591.7334
$ws.Range("K122").Value = 4697.625
$ws.Range("L122").Value = 7775.2002
$ws.Range("M122").Value = -2247.625
$ws.Range("N122").Value = -12675.2002

$ws.Range("H138").Value = 58368.57
$ws.Range("J138").Value = 58368.57
$ws.Range("L138").Value = 58368.57
$ws.Range("N138").Value = -68648.57000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1450377.1
$ws.Range("J131").Value = 1494.3636
$ws.Range("L131").Value = 4483.0908
$ws.Range("N131").Value = -14563.0908

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 8500
$ws.Range("J5").Value = 8500
$ws.Range("L5").Value = 8500
$ws.Range("N5").Value = -8724

$ws.Range("H52").Value = 166684860
$ws.Range("I52").Value = 500005000
$ws.Range("K52").Value = 500005000
$ws.Range("M52").Value = -500004741

$ws.Range("H97").Value = 2343.9092
$ws.Range("I97").Value = 2235.4443
$ws.Range("K97").Value = 2235.4443
$ws.Range("M97").Value = -1739.4443

$ws.Range("H122").Value = 9262302
$ws.Range("I122").Value = 3603178.8
$ws.Range("J122").Value = 15628815
$ws.Range("K122").Value = 10809536.4
$ws.Range("L122").Value = 46886445
$ws.Range("M122").Value = -10807086.4
$ws.Range("N122").Value = -46891345

$ws.Range("H126").Value = 12550
$ws.Range("I126").Value = 15150
$ws.Range("J126").Value = 3450
$ws.Range("K126").Value = 45450
$ws.Range("L126").Value = 10350
$ws.Range("M126").Value = -42980
$ws.Range("N126").Value = -15290

$ws.Range("H133").Value = 42192.223
$ws.Range("J133").Value = 42192.223
$ws.Range("L133").Value = 42192.223
$ws.Range("N133").Value = -52312.223

$ws.Range("H135").Value = 56504.445
$ws.Range("J135").Value = 56504.445
$ws.Range("L135").Value = 56504.445
$ws.Range("N135").Value = -66644.44500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 866.25
$ws.Range("I16").Value = 846
$ws.Range("K16").Value = 846
$ws.Range("M16").Value = -676

$ws.Range("H46").Value = 1020.4
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1020.4
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 1020.4
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -1396.4

$ws.Range("H48").Value = 5007500
$ws.Range("I48").Value = 10000000
$ws.Range("J48").Value = 15000
$ws.Range("K48").Value = 10000000
$ws.Range("L48").Value = 15000
$ws.Range("M48").Value = -9999339
$ws.Range("N48").Value = -16322

$ws.Range("H51").Value = 10250
$ws.Range("J51").Value = 10250
$ws.Range("L51").Value = 10250
$ws.Range("N51").Value = -11206

$ws.Range("H53").Value = 11666.667
$ws.Range("I53").Value = 7000
$ws.Range("J53").Value = 12600
$ws.Range("K53").Value = 7000
$ws.Range("L53").Value = 12600
$ws.Range("M53").Value = -6482
$ws.Range("N53").Value = -13636

$ws.Range("H93").Value = 1709.2727
$ws.Range("I93").Value = 1422.4445
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 1422.4445
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -174.4445000000001
$ws.Range("N93").Value = -5496

$ws.Range("H122").Value = 4292345
$ws.Range("I122").Value = 7946505.5
$ws.Range("J122").Value = 1003600
$ws.Range("K122").Value = 23839516.5
$ws.Range("L122").Value = 3010800
$ws.Range("M122").Value = -23837066.5
$ws.Range("N122").Value = -3015700

$ws.Range("H136").Value = 8708.370999999999
$ws.Range("I136").Value = 5581.9287
$ws.Range("J136").Value = 21214.143
$ws.Range("K136").Value = 16745.7861
$ws.Range("L136").Value = 63642.429
$ws.Range("M136").Value = -14195.7861
$ws.Range("N136").Value = -68742.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 7462.5
$ws.Range("J54").Value = 7462.5
$ws.Range("L54").Value = 7462.5
$ws.Range("N54").Value = -8502.5

$ws.Range("H70").Value = 29095
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 29095
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H81").Value = 1611.2222
$ws.Range("J81").Value = 2466.6667
$ws.Range("L81").Value = 4933.3334
$ws.Range("N81").Value = -7055.3334

$ws.Range("H84").Value = 1611.2222
$ws.Range("J84").Value = 2466.6667
$ws.Range("L84").Value = 24666.667
$ws.Range("N84").Value = -35274.667

$ws.Range("H136").Value = 2351.194
$ws.Range("I136").Value = 2091.875
$ws.Range("J136").Value = 2735.3704
$ws.Range("K136").Value = 6275.625
$ws.Range("L136").Value = 8206.111199999999
$ws.Range("M136").Value = -3725.625
$ws.Range("N136").Value = -13306.1112
